$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 with the "Troll" enemy entry
$ws.Range("A16").Value = "Troll"
$ws.Range("B16").Value = 23
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 6

# Update selection to match target state
$ws.Range("J13").Select() | Out-Null
